$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" column (G) values - replaces former Strike# values
$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 0
    12 = 2
    13 = 0
    14 = 0
    15 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 3
    23 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

$wb.Save()
